$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1128055824"
$ws.Range("D16").Value = "JONATHAN TERAN TORRES"
$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 30284

$ws.Range("C17").Value = "1128055824"
$ws.Range("D17").Value = "JONATHAN TERAN TORRES"
$ws.Range("E17").Value = "2112"
$ws.Range("F17").Value = 36341

$ws.Range("C18").Value = "1128055824"
$ws.Range("D18").Value = "JONATHAN TERAN TORRES"
$ws.Range("E18").Value = "2111"
$ws.Range("F18").Value = 36341

$ws.Range("C19").Value = "1128055824"
$ws.Range("D19").Value = "JONATHAN TERAN TORRES"
$ws.Range("E19").Value = "2110"
$ws.Range("F19").Value = 36341

$ws.Range("C20").Value = "1128055824"
$ws.Range("D20").Value = "JONATHAN TERAN TORRES"
$ws.Range("E20").Value = "2109"
$ws.Range("F20").Value = 36341

$ws.Range("C21").Value = "1128055824"
$ws.Range("D21").Value = "JONATHAN TERAN TORRES"
$ws.Range("E21").Value = "2108"
$ws.Range("F21").Value = 36341

$ws.Range("C22").Value = "1128055824"
$ws.Range("D22").Value = "JONATHAN TERAN TORRES"
$ws.Range("E22").Value = "2107"
$ws.Range("F22").Value = 36341

$ws.Range("C23").Value = "92500773"
$ws.Range("D23").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E23").Value = "2201"
$ws.Range("F23").Value = 30284

$ws.Range("C24").Value = "92500773"
$ws.Range("D24").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E24").Value = "2112"
$ws.Range("F24").Value = 36341

$ws.Range("C25").Value = "92500773"
$ws.Range("D25").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E25").Value = "2111"
$ws.Range("F25").Value = 36341

$ws.Range("C26").Value = "92500773"
$ws.Range("D26").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E26").Value = "2110"
$ws.Range("F26").Value = 36341

$ws.Range("C27").Value = "92500773"
$ws.Range("D27").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E27").Value = "2109"
$ws.Range("F27").Value = 36341

$ws.Range("C28").Value = "92500773"
$ws.Range("D28").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E28").Value = "2108"
$ws.Range("F28").Value = 36341

$ws.Range("C29").Value = "92500773"
$ws.Range("D29").Value = "ALFONSO RAFAEL TERAN MONTES"
$ws.Range("E29").Value = "2107"
$ws.Range("F29").Value = 36341

$ws.Range("C30").Value = "1002322447"
$ws.Range("D30").Value = "LEONARDO ENRIQUE MARIMON SANDOVAL"
$ws.Range("E30").Value = "2107"
$ws.Range("F30").Value = 36341

$ws.Range("C31").Value = "1002257296"
$ws.Range("D31").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E31").Value = "2201"
$ws.Range("F31").Value = 30284

$ws.Range("C32").Value = "1002257296"
$ws.Range("D32").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E32").Value = "2112"
$ws.Range("F32").Value = 36341

$ws.Range("C33").Value = "1002257296"
$ws.Range("D33").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E33").Value = "2111"
$ws.Range("F33").Value = 36341

$ws.Range("C34").Value = "1002257296"
$ws.Range("D34").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E34").Value = "2110"
$ws.Range("F34").Value = 36341

$ws.Range("C35").Value = "1002257296"
$ws.Range("D35").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E35").Value = "2109"
$ws.Range("F35").Value = 36341

$ws.Range("C36").Value = "1002257296"
$ws.Range("D36").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E36").Value = "2108"
$ws.Range("F36").Value = 36341

$ws.Range("C37").Value = "1002257296"
$ws.Range("D37").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E37").Value = "2107"
$ws.Range("F37").Value = 36341

$ws.Range("C38").Value = "1002257296"
$ws.Range("D38").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E38").Value = "2106"
$ws.Range("F38").Value = 36341

$ws.Range("C39").Value = "1002257296"
$ws.Range("D39").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E39").Value = "2105"
$ws.Range("F39").Value = 36341

$ws.Range("C40").Value = "1002257296"
$ws.Range("D40").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E40").Value = "2104"
$ws.Range("F40").Value = 36341

$ws.Range("C41").Value = "1002257296"
$ws.Range("D41").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E41").Value = "2103"
$ws.Range("F41").Value = 36341

$ws.Range("C42").Value = "1002257296"
$ws.Range("D42").Value = "MAICKOL ANDRES POLO SILVA"
$ws.Range("E42").Value = "2102"
$ws.Range("F42").Value = 36341
